$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.621.29'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.56%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.601.53'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.45%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.52'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.10%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.515'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.88%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.23%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '27.18'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +4.23%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +1.20%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +1.33%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.50%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.829.87'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.29%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.600.39'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.91%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.542'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +3.83%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '29.622.49'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.46%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.75'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.30%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.73'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +2.50%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '241.27'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.94%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.65'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +3.01%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0695'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.77%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.22%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.00'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.37%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.30'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +1.51%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.45%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '155.00'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.83%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.39'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +1.65%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.71%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.41'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.03%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.18%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0479'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +2.35%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.08%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.62%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.18'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +3.71%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.420.52'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.03%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +2.93%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +4.01%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -1.96%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.09%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +3.32%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.546'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +3.43%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '56.02'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +5.19%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.46%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +4.66%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.815'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +3.40%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.19%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +17.04%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '65.91'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +2.28%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.32'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.20%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.740.57'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +1.30%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '86.25'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.07%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +3.79%  '
